$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet2" to "Sheet1"
$ws.Name = "Sheet1"

# Clear the "#" index cell in A2 (row number column no longer used)
$ws.Range("A2").ClearContents()

# Recreate the mailto hyperlink on C2 so it keeps pointing at icists@icists.org
# while showing a friendly display label (the cell text itself is updated below
# to list both recipients).
$ws.Range("C2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:icists@icists.org", "", "", "icists@icists.org")
$ws.Range("C2").Style = "하이퍼링크"

# Update the mail address cell (C2) to include the additional recipient
$ws.Range("C2").Value = "icists@icists.org, media@icists.org"

# Move the active selection to C3, matching the saved selection state
$ws.Range("C3").Select()
